$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append rows 27-31 ----
$pir = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-01-30", "16:01:30", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:01:30", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:01:34", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:01:39", "16:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "16:01:44", "16:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 27
# Force the Date column to be stored as plain text so values like
# "2026-01-30" are not auto-converted into date serial numbers.
$pir.Range("A" + $startRow + ":A" + ($startRow + $pirRows.Count - 1)).NumberFormat = "@"

for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]
    $pir.Cells.Item($r, 1).Value = $row[0]
    $pir.Cells.Item($r, 2).Value = $row[1]
    $pir.Cells.Item($r, 3).Value = $row[2]
    $pir.Cells.Item($r, 4).Value = $row[3]
    $pir.Cells.Item($r, 5).Value = $row[4]
    $pir.Cells.Item($r, 6).Value = $row[5]
}

# ---- Humidity sheet: append rows 21-23 ----
$hum = $wb.Worksheets.Item("Humidity")

$humRows = @(
    @("2026-01-30", "16:01:30", "16:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "16:01:35", "16:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "16:01:45", "16:00", "Bathroom", "87.6%", "Active")
)

$startRow = 21
# Force Date (A) and Value (E, a percentage-as-text reading) columns to be
# stored as plain text, matching the original rows' literal text content.
$hum.Range("A" + $startRow + ":A" + ($startRow + $humRows.Count - 1)).NumberFormat = "@"
$hum.Range("E" + $startRow + ":E" + ($startRow + $humRows.Count - 1)).NumberFormat = "@"

for ($i = 0; $i -lt $humRows.Count; $i++) {
    $r = $startRow + $i
    $row = $humRows[$i]
    $hum.Cells.Item($r, 1).Value = $row[0]
    $hum.Cells.Item($r, 2).Value = $row[1]
    $hum.Cells.Item($r, 3).Value = $row[2]
    $hum.Cells.Item($r, 4).Value = $row[3]
    $hum.Cells.Item($r, 5).Value = $row[4]
    $hum.Cells.Item($r, 6).Value = $row[5]
}
